$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row before the existing row 138, shifting all
# subsequent rows (138-244) down by one (they become 139-245).
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new weekly record.
$ws.Range("A138").Value = 5
$ws.Range("B138").Value = "Macroferia Regional de Talca"
$ws.Range("C138").Value = "Maule"
$ws.Range("D138").Value = 44574
$ws.Range("D138").NumberFormat = $ws.Range("D139").NumberFormat
$ws.Range("E138").Value = 7
$ws.Range("F138").Value = 100114014
$ws.Range("G138").Value = "Betarraga"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 500
$ws.Range("L138").Value = 500
$ws.Range("M138").Value = 500
$ws.Range("N138").Value = "$/paquete 5 unidades"
$ws.Range("O138").Value = "Región del Maule"
$ws.Range("P138").Value = 100
$ws.Range("Q138").Value = 5
$ws.Range("R138").Value = "Hortaliza"
